$d = $word.ActiveDocument

# --- 1. Locate the day portion ("25") inside the "Fecha: 25/02/16" line ---
$dateRange = $d.Content
$found = $dateRange.Find.Execute("25/02/16", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the date '25/02/16' to edit"
}

$dayStart = $dateRange.Start
$dayEnd   = $dayStart + 2          # "25" is two characters long

# --- 2. Drop a throw-away bookmark exactly on the existing run boundary that
#        precedes "25" (i.e. right after "Fecha:"). This stops the engine
#        from re-merging that earlier run into the one we are about to edit,
#        so "Fecha:" stays in its own run exactly as before. ---
$barrierRange = $d.Range($dayStart - 1, $dayStart - 1)
$d.Bookmarks.Add("zzz_barrier_tmp", $barrierRange)

# --- 3. Move (or create) the "_GoBack" bookmark to sit right after the new
#        day value -- this is where Word leaves it after the user's last
#        keystroke. Adding it here also splits the run at that position,
#        producing a fresh run that will hold the untouched "/02/16" tail. ---
$goBackRange = $d.Range($dayEnd, $dayEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# --- 4. Replace "25" with "29" ---
$dayRange = $d.Range($dayStart, $dayEnd)
$dayRange.Text = "29"

# --- 5. Remove the temporary barrier bookmark; it has done its job ---
$d.Bookmarks("zzz_barrier_tmp").Delete()

# --- 6. Touch the trailing "/02/16" run so the engine recomputes its
#        serialization (drops the now-unnecessary xml:space="preserve" it
#        would otherwise inherit verbatim from the original, pre-split run)
#        while keeping its original character formatting intact. ---
$tailStart = $dayEnd
$tailEnd   = $tailStart + 6        # "/02/16" is six characters long
$tailRange = $d.Range($tailStart, $tailEnd)
$tailRange.Text = "zzzzzz"
$tailRange2 = $d.Range($tailStart, $tailStart + 6)
$tailRange2.Text = "/02/16"
